# config for new simulations
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# realistic_candidate_capacities_for_future -> FALSE
$ws.Range("B20").Value = $false

# dummy_capacity -> 1
$ws.Range("B22").Value = 1

# move the active selection to B21
$ws.Range("B21").Select()
